# Applies the "add outliers handling for continuous metrics" edit.
#
# Summary of the change (see the unified diff used to derive this script):
#  - sheet "_all": row5 "Unique ids" count 6 -> 9; the "conversion to
#    purchase" and "CTR" comparison rows (11/12) get new recomputed
#    percentages/pvalues; the corresponding corrected-pvalue rows (17-20)
#    get new pvalues/flags.
#  - sheet "android": row5 "Unique ids" count 6 -> 9 (no other changes).
#  - sheet "ios": row5 "Unique ids" count 6 -> 9; the ARPU row (11) gets
#    real percentages, and a brand new "ARPU ... - quantile 0.98" metric
#    is inserted as row 12 (with its own corrected-pvalue rows 19-20),
#    pushing the old "metric" header block from row 15 down to row 16.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "_all"
# ---------------------------------------------------------------------
$all = $wb.Worksheets.Item("_all")

$all.Range("A5").Value = "'Unique ids: 9"

$all.Range("B11").Value = "'50.00%"
$all.Range("C11").Value = "'40.00%"
$all.Range("E11").Value = "'-20.00%"
$all.Range("F11").Value = "'+100.00%"
$all.Range("G11").Value = 0.7641771556220945
$all.Range("H11").Value = 0.4413427238396865

$all.Range("B12").Value = "'76.00%"
$all.Range("C12").Value = "'113.64%"
$all.Range("D12").Value = "'104.76%"
$all.Range("E12").Value = "'+49.52%"
$all.Range("F12").Value = "'+37.84%"
$all.Range("G12").Value = 0.0000005551375349886539
$all.Range("H12").Value = 0.005393233202449595
$all.Range("I12").Value = $true
$all.Range("J12").Value = $true

$all.Range("C17").Value = 0.7641771556220945
$all.Range("F17").Value = 0.7641771556220945

$all.Range("C18").Value = 0.2206713619198432
$all.Range("F18").Value = 0.4413427238396865

$all.Range("C19").Value = 0.0000001387843837471635
$all.Range("E19").Value = $true
$all.Range("F19").Value = 0.0000005551375349886539
$all.Range("G19").Value = 1

$all.Range("C20").Value = 0.001797744400816532
$all.Range("E20").Value = $true
$all.Range("F20").Value = 0.005393233202449595
$all.Range("G20").Value = 1

# ---------------------------------------------------------------------
# Sheet "android"
# ---------------------------------------------------------------------
$android = $wb.Worksheets.Item("android")
$android.Range("A5").Value = "'Unique ids: 9"

# ---------------------------------------------------------------------
# Sheet "ios"
# ---------------------------------------------------------------------
$ios = $wb.Worksheets.Item("ios")
$ios.Range("A5").Value = "'Unique ids: 9"

# Insert the new "quantile" metric row right after the existing ARPU row
# (row 11), pushing the old rows 15-17 down to 16-18.
$ios.Rows.Item(12).Insert()

# Fill in the real percentages on the (now complete) ARPU row.
$ios.Range("B11").Value = "'110.000"
$ios.Range("C11").Value = "'210.000"
$ios.Range("D11").Value = "'145.000"
$ios.Range("E11").Value = "'+90.91%"
$ios.Range("F11").Value = "'+31.82%"

# New row 12: "ARPU (USD after fee) - quantile 0.98"
$ios.Range("A11").Copy()
$ios.Range("A12").PasteSpecial(-4122)
$ios.Range("A12").Value = "ARPU (USD after fee) - quantile 0.98"
$ios.Range("B12").Value = "'110.000"
$ios.Range("C12").Value = "'210.000"
$ios.Range("D12").Value = "'15.000"
$ios.Range("E12").Value = "'+90.91%"
$ios.Range("F12").Value = "'-86.36%"

# New row 18 ("3-1" row for the original ARPU metric, format copied from
# the sibling "3-1" row pattern used elsewhere in the sheet).
$ios.Range("A17:B17").Copy()
$ios.Range("A18").PasteSpecial(-4122)
$ios.Range("B18").Value = "'3-1"
$ios.Range("D18").Value = 0.05

# New rows 19-20: corrected-pvalue block for the new quantile metric.
$ios.Range("A17:B18").Copy()
$ios.Range("A19").PasteSpecial(-4122)
$ios.Range("A19").Value = "ARPU (USD after fee) - quantile 0.98"
$ios.Range("B19").Value = "'2-1"
$ios.Range("D19").Value = 0.05

$ios.Range("B20").Value = "'3-1"
$ios.Range("D20").Value = 0.05

$ios.Range("A19:A20").Merge()
